$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 11) of data for "2021年" below the existing last row (row 10, "2020年"),
# following the same column layout as the rest of the table.
$ws.Range("A11").Value = "2021年"
$ws.Range("B11").Value = 15050.75
$ws.Range("C11").Value = 3183.45
$ws.Range("D11").Value = 1706.68
$ws.Range("F11").Value = 7993.46
$ws.Range("G11").Value = 30610.42
$ws.Range("H11").Value = 4585.39
$ws.Range("I11").Value = 11737.73
$ws.Range("J11").Value = 1176.1
$ws.Range("K11").Value = 550415.87
$ws.Range("L11").Value = 2223.35
$ws.Range("M11").Value = 764.6
$ws.Range("N11").Value = 1815.92
$ws.Range("O11").Value = 2508.6
$ws.Range("P11").Value = 20451.24
$ws.Range("Q11").Value = 2310.45
$ws.Range("R11").Value = 738.59
$ws.Range("S11").Value = 6294.33
$ws.Range("T11").Value = 7235.86
$ws.Range("U11").Value = 41524.18
$ws.Range("V11").Value = 2452.71
$ws.Range("W11").Value = 38414.58
$ws.Range("X11").Value = 4062.24
$ws.Range("Y11").Value = 71723.34
$ws.Range("Z11").Value = 36963.04
$ws.Range("AA11").Value = 1685.93
$ws.Range("AB11").Value = 23784.65
$ws.Range("AC11").Value = 9733.84
$ws.Range("AD11").Value = 5668.95
$ws.Range("AE11").Value = 3155.17
$ws.Range("AF11").Value = 70691.66
$ws.Range("AG11").Value = 18397.76
$ws.Range("AH11").Value = 6001.19
$ws.Range("AI11").Value = 5924.48
$ws.Range("AJ11").Value = 1019.57
$ws.Range("AK11").Value = 10586.32
$ws.Range("AL11").Value = 14938.48
$ws.Range("AM11").Value = 13843.35
$ws.Range("AN11").Value = 718.34
$ws.Range("AO11").Value = 5837.56
$ws.Range("AP11").Value = 37293.62
$ws.Range("AQ11").Value = 5603.97

# Column E is left blank for this row, matching the blank cell already used
# in column E of the previous row (E10).

# Copy the formatting of the label cell in the row above (A10, which carries
# the bold/centered/bordered style used for all year labels in column A)
# onto the new label cell A11, without disturbing the value we just set.
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
